$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.967707753181458
$ws.Range("B1").Value = 2.120103359222412
$ws.Range("C1").Value = 2.171495914459229
$ws.Range("D1").Value = 2.643497943878174
$ws.Range("E1").Value = 3.772227764129639
